$d = $word.ActiveDocument

# Target: Table 1, row 75 ("monoterpenoid emissions" / "MT" / <empty units cell>),
# the empty "units" cell (column 3). It needs "µg m-2s-1" with the two
# exponents set as italic superscripts, matching the look of
# "µg m<sup>-2</sup> s<sup>-1</sup>".
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(75, 3)

# Insert the full plain run first (InsertAfter on the Cell.Range keeps the
# target position anchored correctly, unlike building fresh Range objects
# before any text exists in the cell).
$cell.Range.InsertAfter("µg m-2s-1")

# Re-fetch the cell range now that it holds text, and locate the text
# span (End points one past the last character, before the cell mark).
$cellRange = $tbl.Cell(75, 3).Range
$cellStart = $cellRange.Start
$cellTextEnd = $cellRange.End - 1

# Segment offsets within "µg m-2s-1":
#   "µg m" -> 0..4   "-2" -> 4..6   "s" -> 6..7   "-1" -> 7..9
$seg1 = $d.Range($cellStart + 0, $cellStart + 4)   # "µg m"
$seg2 = $d.Range($cellStart + 4, $cellStart + 6)   # "-2"
$seg3 = $d.Range($cellStart + 6, $cellStart + 7)   # "s"
$seg4 = $d.Range($cellStart + 7, $cellTextEnd)     # "-1"

$seg1.Font.Name = "Calibri"
$seg1.Font.Italic = $true
$seg1.Font.Color = 0

$seg2.Font.Name = "Calibri"
$seg2.Font.Italic = $true
$seg2.Font.Color = 0
$seg2.Font.Superscript = $true
$seg2.Font.Size = 12

$seg3.Font.Name = "Calibri"
$seg3.Font.Italic = $true
$seg3.Font.Color = 0

$seg4.Font.Name = "Calibri"
$seg4.Font.Italic = $true
$seg4.Font.Color = 0
$seg4.Font.Superscript = $true
$seg4.Font.Size = 12
